{"js": "// Add two new character paragraphs (\"Lion\", \"Tiger\") right after the\n// existing \"Crocodile\" paragraph, before the end of the document body.\nconst body = context.document.body;\n\n// Insert \"Lion\" at the end of the body (after \"Crocodile\"), then \"Tiger\"\n// at the new end (after \"Lion\") -- each insertParagraph(\"End\") call lands\n// the new paragraph just before the body's closing section properties,\n// matching Word's own behavior for appending paragraphs.\nbody.insertParagraph(\"Lion\", Word.InsertLocation.end);\nbody.insertParagraph(\"Tiger\", Word.InsertLocation.end);\n\nawait context.sync();\n", "ps1": "# Add two new character paragraphs (\"Lion\", \"Tiger\") right after the\n# existing \"Crocodile\" paragraph, at the end of the document body.\n$d = $word.ActiveDocument\n\n# Paragraphs.Add() appends a brand-new paragraph at the end of the story;\n# calling it twice in sequence lands \"Lion\" then \"Tiger\" directly after\n# \"Crocodile\", just before the section break -- matching Word's own\n# behavior for appending paragraphs at the end of a document.\n$p1 = $d.Paragraphs.Add()\n$p1.Range.Text = \"Lion\"\n\n$p2 = $d.Paragraphs.Add()\n$p2.Range.Text = \"Tiger\"\n"}
